$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking text values (prices / percentages) from being
# auto-coerced to numbers by formatting the affected range as Text first.
$numRange = $ws.Range("D2:E51")
$numRange.NumberFormat = "@"

$ws.Range("D2").Value = "331.84"
$ws.Range("E2").Value = "0.50%"
$ws.Range("D3").Value = "45.41"
$ws.Range("E3").Value = "2.80%"
$ws.Range("D4").Value = "5.612"
$ws.Range("E4").Value = "2.36%"
$ws.Range("D5").Value = "0.08352"
$ws.Range("E5").Value = "4.38%"
$ws.Range("D6").Value = "2.086"
$ws.Range("E6").Value = "5.67%"
$ws.Range("D7").Value = "0.9660"
$ws.Range("E7").Value = "1.49%"
$ws.Range("D8").Value = "2.553"
$ws.Range("E8").Value = "-1.02%"
$ws.Range("E9").Value = "4.15%"
$ws.Range("D10").Value = "0.1924"
$ws.Range("E10").Value = "1.10%"
$ws.Range("D11").Value = "10.30"
$ws.Range("E11").Value = "-4.36%"
$ws.Range("D12").Value = "0.09852"
$ws.Range("E12").Value = "-1.16%"
$ws.Range("D13").Value = "0.04619"
$ws.Range("E13").Value = "-3.50%"
$ws.Range("E14").Value = "-0.18%"
$ws.Range("D15").Value = "0.001280"
$ws.Range("E15").Value = "0.17%"
$ws.Range("D16").Value = "0.006119"
$ws.Range("E16").Value = "3.09%"
$ws.Range("D17").Value = "3.376"
$ws.Range("D18").Value = "4.442"
$ws.Range("E18").Value = "1.27%"
$ws.Range("D19").Value = "0.3349"
$ws.Range("E19").Value = "-4.03%"
$ws.Range("D20").Value = "0.1394"
$ws.Range("E20").Value = "-1.79%"
$ws.Range("D21").Value = "0.2881"
$ws.Range("E21").Value = "11.50%"
$ws.Range("D22").Value = "0.04186"
$ws.Range("E22").Value = "2.53%"
$ws.Range("E23").Value = "3.57%"
$ws.Range("D24").Value = "0.004558"
$ws.Range("E24").Value = "4.77%"
$ws.Range("E25").Value = "8.73%"
$ws.Range("D26").Value = "0.0003752"
$ws.Range("E26").Value = "0.27%"
$ws.Range("D38").Value = "0.02707"
$ws.Range("E38").Value = "4.47%"
$ws.Range("D39").Value = "0.05759"
$ws.Range("E39").Value = "0.40%"
$ws.Range("D40").Value = "0.007899"
$ws.Range("E40").Value = "4.43%"
$ws.Range("E41").Value = "2.34%"
$ws.Range("D42").Value = "0.007272"
$ws.Range("E42").Value = "-1.14%"
$ws.Range("E43").Value = "0.60%"
$ws.Range("D44").Value = "0.009114"
$ws.Range("E44").Value = "3.48%"
$ws.Range("D45").Value = "0.3545"
$ws.Range("D46").Value = "0.00007134"
$ws.Range("E46").Value = "0.18%"
$ws.Range("E47").Value = "0.38%"
$ws.Range("D48").Value = "0.0005819"
$ws.Range("E48").Value = "0.38%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "0.003509"
$ws.Range("E49").Value = "-0.53%"
$ws.Range("B50").Value = "BOLO"
$ws.Range("C50").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D50").Value = "0.003495"
$ws.Range("E50").Value = "-1.66%"
$ws.Range("D51").Value = "0.00002107"
$ws.Range("E51").Value = "0.38%"

# Restore the default (unstyled) cell style now that values are committed,
# so no stray number-format styling is left behind on the data cells.
$numRange.Style = "Normal"
